# New weekly price record for "Brócoli" at Feria Lagunitas de Puerto Montt.
# The record belongs chronologically right before the current row 313
# (date 44659), so insert a new row there and push the existing rows
# (313-338) down to (314-339).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(313).Insert()

$ws.Range("A313").Value = 4
$ws.Range("B313").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C313").Value = "Los Lagos"
$ws.Range("D313").Value = 44714
$ws.Range("E313").Value = 10
$ws.Range("F313").Value = 100112023
$ws.Range("G313").Value = "Brócoli"
$ws.Range("H313").Value = "Sin especificar"
$ws.Range("I313").Value = "Primera"
$ws.Range("J313").Value = 500
$ws.Range("K313").Value = 1400
$ws.Range("L313").Value = 1500
$ws.Range("M313").Value = 1450
$ws.Range("N313").Value = "$/unidad"
$ws.Range("O313").Value = "Región Metropolitana"
$ws.Range("P313").Value = 1450
$ws.Range("Q313").Value = 1
$ws.Range("R313").Value = "Hortaliza"
